# Add a new "CATEGORIES" section (model/controller/tests/view) to the API
# documentation worksheet, mirroring the layout of the existing sections
# (AGENTS, COMPANIES, CONTACTS, INCIDENTS).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant used with PasteSpecial below
$xlPasteFormats = -4122

# --- Row 42: section title "CATEGORIES" (same look as the other section
#     headers in A1/A11/A23/A31 -> style index 2) ---
$ws.Range("A11").Copy()
$ws.Range("A42").PasteSpecial($xlPasteFormats)
$ws.Range("A42").Value = "CATEGORIES"
$ws.Rows.Item(42).RowHeight = 18

# --- Row 43: first data row of the section ("top" row style, same as
#     A4/B4/C4/D4, A12/B12/C12/D12, etc. -> style indices 4,5,5,6) ---
$ws.Range("A4").Copy()
$ws.Range("A43").PasteSpecial($xlPasteFormats)
$ws.Range("A43").Value = "GET"

$ws.Range("B4").Copy()
$ws.Range("B43").PasteSpecial($xlPasteFormats)
$ws.Range("B43").Value = "/categories"

$ws.Range("C4").Copy()
$ws.Range("C43").PasteSpecial($xlPasteFormats)
$ws.Range("C43").Value = "Get complete list of categories"

$ws.Range("D4").Copy()
$ws.Range("D43").PasteSpecial($xlPasteFormats)

# --- Row 44: last data row of the section. It needs a bottom border like
#     the other "last row" styles (A21/A29/A40 etc, borders 6/7/8) but
#     keeps the blue GET fill instead of the green/red ones used
#     elsewhere, so it is a brand-new style combination. We get there by
#     copying the border/font from an existing "last row" style and then
#     only touching the Interior color, so the existing font (fontId=1)
#     is reused rather than a new one being created. ---
$ws.Range("A21").Copy()
$ws.Range("A44").PasteSpecial($xlPasteFormats)
$ws.Range("A44").Interior.Color = 16737843
$ws.Range("A44").Value = "GET"

$ws.Range("B21").Copy()
$ws.Range("B44").PasteSpecial($xlPasteFormats)
$ws.Range("B44").Interior.Color = 16737843
$ws.Range("B44").Value = "/categories/{id}"

$ws.Range("C21").Copy()
$ws.Range("C44").PasteSpecial($xlPasteFormats)
$ws.Range("C44").Interior.Color = 16737843
$ws.Range("C44").Value = "Get information for a specific category"

$ws.Range("D21").Copy()
$ws.Range("D44").PasteSpecial($xlPasteFormats)
$ws.Range("D44").Interior.Color = 16737843

# --- View state: scroll down a bit and select the (now empty) row right
#     after the new section, matching where the author's cursor ended up ---
$ws.Range("A46").Select()
